$wb = $excel.ActiveWorkbook

# "clusters" sheet: shorten the cluster name labels in column A
$ws = $wb.Worksheets.Item("clusters")
$ws.Range("A1").Value = "conduction"
$ws.Range("A2").Value = "zaakonline"

# Make "clusters" the active/selected sheet (was "environments")
$ws.Activate()
[void]$ws.Range("H6").Select()
